$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("2023-08-16" week): D6 (Vendas de Chai pré-fabricado) and E6
# (Engajamento em redes sociais) were pasted in from another source as
# clock-time-looking text ("4:36" / "17:05") instead of the previous
# plain numbers (436 / 1705). Write them as text values, matching the
# new shared-string entries.
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "17:05"
